$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63, pushing existing rows 63+ down by one.
$ws.Rows.Item(63).Insert()

# Populate the new row 63 with the new weekly record.
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44706
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = "Fruta"
$ws.Cells.Item(63, 7).Value = 100108
$ws.Cells.Item(63, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(63, 9).Value = 100108007
$ws.Cells.Item(63, 10).Value = "Coco"
$ws.Cells.Item(63, 11).Value = "Sin especificar"
$ws.Cells.Item(63, 12).Value = "Primera"
$ws.Cells.Item(63, 13).Value = 25
$ws.Cells.Item(63, 14).Value = 30000
$ws.Cells.Item(63, 15).Value = 30000
$ws.Cells.Item(63, 16).Value = 30000
$ws.Cells.Item(63, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(63, 18).Value = "Perú"
$ws.Cells.Item(63, 19).Value = 1500
$ws.Cells.Item(63, 20).Value = 20
